# إضافة حدث جديد في Card20
# Row 13: columns B..K were blank placeholders -> now explicitly "nan"
# Row 14 (new): a new service event row

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card20")

# --- Fill row 13, columns B..K with literal "nan" text ---
$ws.Cells.Item(13, 2).Value  = "nan"   # B13
$ws.Cells.Item(13, 3).Value  = "nan"   # C13
$ws.Cells.Item(13, 4).Value  = "nan"   # D13
$ws.Cells.Item(13, 5).Value  = "nan"   # E13
$ws.Cells.Item(13, 6).Value  = "nan"   # F13
$ws.Cells.Item(13, 7).Value  = "nan"   # G13
$ws.Cells.Item(13, 8).Value  = "nan"   # H13
$ws.Cells.Item(13, 9).Value  = "nan"   # I13
$ws.Cells.Item(13, 10).Value = "nan"   # J13
$ws.Cells.Item(13, 11).Value = "nan"   # K13

# --- Add new row 14 with the new service event ---
# Copy A13 (already the text "20") down into A14 so the new cell keeps
# text storage instead of being re-inferred as a number.
$ws.Range("A13").Copy($ws.Range("A14"))                                   # A14 = "20"
$ws.Cells.Item(14, 12).Value = "10\8\2024"                                # L14
$ws.Cells.Item(14, 14).Value = "تم تشحيم المكنه بالكامل +عمل صيانه"        # N14
$ws.Cells.Item(14, 15).Value = "تيم العمل"                                # O14
